# Patton's Best - Events sheet: add e062 and e105 crew-replacement events
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new event rows, keeping the sheet sorted by event id:
#   ... e061a  [+e062]  e100 ... e104  [+e105]  e501 ...
$ws.Rows.Item(94).Insert()
$ws.Rows.Item(101).Insert()

$ws.Range("A94").Value = "e062"
$ws.Range("A101").Value = "e105"

$ws.Range("B101").Value = "<Bold>e105 Crew Replacement - Knocked Out Tank</Bold> `r`n<InlineUIContainer><Button Content='r7.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `r`n<InlineUIContainer><Button Content='r19.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>     `r`n<LineBreak/><LineBreak/>`r`nWounded crewmen are replaced at beginning of day when the tank is knocked out. Click image to  continue to assign new crew ratings.`r`n<LineBreak/><LineBreak/>"
$ws.Range("B94").Value = "<Bold>e062 Crew Replacement - Battle Ends</Bold> `r`n<InlineUIContainer><Button Content='r7.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `r`n<InlineUIContainer><Button Content='r19.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   `r`n<LineBreak/><LineBreak/>`r`nWounded crewmen are replaced when battle ends or the tank withdraws. Replacing crewmen takes 30 minutes and is marked off the After Action Report `r`n<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.`r`n <LineBreak/><LineBreak/>`r`nClick image to assign replacement crew ratings.`r`n<LineBreak/><LineBreak/>`r`n                                                 <InlineUIContainer><Image Name='CarryingMan' Height='200' Width='80'></Image></InlineUIContainer>"

$ws.Range("B94").RowHeight = 150
$ws.Range("B101").RowHeight = 90

$ws.Range("B92").Select()
